$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before current row 6 (old row6->row8, old row7->row9)
$ws.Rows.Item(6).Resize(2).Insert()

# Enter the "name" values first (this establishes the shared-string order
# seen in the target file: Test 6, ytfd, PRJ-03, Test Wrong ID)
$ws.Cells.Item(6, 2).Value = "Test 6"
$ws.Cells.Item(7, 2).Value = "ytfd"
$ws.Cells.Item(6, 4).Value = "PRJ-03"
$ws.Cells.Item(7, 4).Value = "PRJ-03"
$ws.Cells.Item(10, 2).Value = "Test Wrong ID"

# Row 6: 1006, Test 6, 512, PRJ-03
$ws.Cells.Item(6, 1).Value = 1006
$ws.Cells.Item(6, 3).Value = 512

# Row 7: 1004, ytfd, 167, PRJ-03
$ws.Cells.Item(7, 1).Value = 1004
$ws.Cells.Item(7, 3).Value = 167

# Row 10: 101, Test Wrong ID, 54, PRJ-01
$ws.Cells.Item(10, 1).Value = 101
$ws.Cells.Item(10, 3).Value = 54
$ws.Cells.Item(10, 4).Value = "PRJ-01"

# Set custom (best-fit) width on column B to fit the long string in B8
# (131 is the closest input this engine's width quantization maps back to
# the original author's stored width of 131.85546875)
$ws.Columns.Item(2).ColumnWidth = 131

# Update selection to A11
$ws.Range("A11").Select()
